$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118, shifting rows 118:138 down to 119:139
$ws.Rows("118:118").Insert()

# Populate the new row 118 with the inserted record
$ws.Range("A118").Value = 5
$ws.Range("B118").Value = "Macroferia Regional de Talca"
$ws.Range("C118").Value = "Maule"
$ws.Range("D118").Value = 44776
$ws.Range("E118").Value = 7
$ws.Range("F118").Value = "Fruta"
$ws.Range("G118").Value = 100108
$ws.Range("H118").Value = "Tropicales y subtropicales"
$ws.Range("I118").Value = 100108002
$ws.Range("J118").Value = "Mango"
$ws.Range("K118").Value = "Sin especificar"
$ws.Range("L118").Value = "Primera"
$ws.Range("M118").Value = 228
$ws.Range("N118").Value = 10000
$ws.Range("O118").Value = 10000
$ws.Range("P118").Value = 10000
$ws.Range("Q118").Value = '$/bandeja 4 kilos'
$ws.Range("R118").Value = "Brasil"
$ws.Range("S118").Value = 2500
$ws.Range("T118").Value = 4
